# Updated cryptos list on Thu Jan  4 18:30:59 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "44.051.17"
$ws.Range("E2").Value = "  +2.21%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.275.28"
$ws.Range("E3").Value = "  +2.61%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.32%  "

# Row 5 - BNB
Set-TextCell "D5" "319.22"
$ws.Range("E5").Value = "  +0.88%  "

# Row 6 - Solana
Set-TextCell "D6" "103.37"
$ws.Range("E6").Value = "  +5.26%  "

# Row 7 - XRP
Set-TextCell "D7" "0.588"
$ws.Range("E7").Value = "  +0.96%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.35%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.572"
$ws.Range("E9").Value = "  +2.05%  "

# Row 10 - Avalanche
Set-TextCell "D10" "38.87"
$ws.Range("E10").Value = "  +5.96%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.63%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  +1.32%  "

# Row 13 - TRON
Set-TextCell "D13" "0.107"
$ws.Range("E13").Value = "  +1.76%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextCell "D14" "2.622.92"
$ws.Range("E14").Value = "  +2.49%  "

# Row 15 - Polygon
Set-TextCell "D15" "0.876"
$ws.Range("E15").Value = "  +1.10%  "

# Row 16 - Chainlink
Set-TextCell "D16" "14.56"
$ws.Range("E16").Value = "  +3.19%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "2.277.13"
$ws.Range("E17").Value = "  +2.97%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "43.957.84"
$ws.Range("E18").Value = "  +2.31%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextCell "D19" "14.17"
$ws.Range("E19").Value = "  -3.69%  "

# Row 20 - ShibaInu
Set-TextCell "D20" "0.0₃0999"
$ws.Range("E20").Value = "  +4.08%  "

# Row 21 - Uniswap
Set-TextCell "D21" "6.64"
$ws.Range("E21").Value = "  +3.58%  "

# Row 22 - Litecoin
Set-TextCell "D22" "66.21"
$ws.Range("E22").Value = "  +1.42%  "

# Row 23 - PancakeSwap
$ws.Range("E23").Value = "  +0.30%  "

# Row 24 - BitcoinCash
Set-TextCell "D24" "237.17"
$ws.Range("E24").Value = "  -0.01%  "

# Row 25 - ImmutableX
Set-TextCell "D25" "2.19"
$ws.Range("E25").Value = "  +3.05%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.03%  "

# Row 27 - Cosmos
Set-TextCell "D27" "10.26"
$ws.Range("E27").Value = "  +1.78%  "

# Row 28 - InjectiveProtocol
Set-TextCell "D28" "38.75"
$ws.Range("E28").Value = "  +14.44%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -0.03%  "

# Row 30 - Filecoin
Set-TextCell "D30" "6.56"
$ws.Range("E30").Value = "  +4.63%  "

# Row 31 - Monero
Set-TextCell "D31" "162.28"
$ws.Range("E31").Value = "  +4.68%  "

# Row 32 - EthereumClassic
$ws.Range("E32").Value = "  +0.17%  "

# Row 33 - Hedera
Set-TextCell "D33" "0.0878"
$ws.Range("E33").Value = "  -0.78%  "

# Row 34 - WEMIXToken
Set-TextCell "D34" "2.72"
$ws.Range("E34").Value = "  -2.48%  "

# Rows 35 & 36 swap: LidoDAOToken <-> ARBITRUM
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D35" "2.06"
$ws.Range("E35").Value = "  +2.18%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D36" "3.22"
$ws.Range("E36").Value = "  -0.20%  "

# Row 37 - Stellar
$ws.Range("E37").Value = "  -0.87%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  +5.32%  "

# Row 39 - RenderToken
Set-TextCell "D39" "4.50"
$ws.Range("E39").Value = "  +1.16%  "

# Row 40 - NEARProtocol
Set-TextCell "D40" "3.91"
$ws.Range("E40").Value = "  +5.91%  "

# Rows 41 & 42 swap: Celestia <-> VeChain
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D41" "0.0326"
$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell "D42" "15.38"
$ws.Range("E42").Value = "  +24.99%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  -0.25%  "

# Row 44 - Maker
Set-TextCell "D44" "1.769.16"
$ws.Range("E44").Value = "  -5.54%  "

# Row 45 - Algorand
Set-TextCell "D45" "0.207"
$ws.Range("E45").Value = "  +0.15%  "

# Row 46 - BitcoinSV
Set-TextCell "D46" "85.15"
$ws.Range("E46").Value = "  -4.33%  "

# Row 47 - THORChain
Set-TextCell "D47" "5.36"
$ws.Range("E47").Value = "  -1.60%  "

# Row 48 - FraxShare
Set-TextCell "D48" "8.87"
$ws.Range("E48").Value = "  +2.18%  "

# Row 49 - MultiversX
$ws.Range("E49").Value = "  -1.13%  "

# Row 50 - ordi
Set-TextCell "D50" "74.79"
$ws.Range("E50").Value = "  -2.83%  "

# Row 51 - Aave
Set-TextCell "D51" "104.44"
$ws.Range("E51").Value = "  +3.71%  "
